$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(31.0, "Anny Rectangular Metal Coffee Table In Powder Coating Finish", "₹23,220")
    ,@(32.0, "Lavish Rectangular Metal Coffee Table In Powder Coating Finish", "₹23,220")
    ,@(33.0, "Frazer Rectangular Metal Coffee Table In Powder Coating Finish", "₹25,451")
    ,@(34.0, "Marten Rectangular Metal Coffee Table In Powder Coating Finish", "₹23,853")
    ,@(35.0, "Alix Rectangular Metal Coffee Table In Stainless Steel Finish", "₹19,763")
    ,@(36.0, "Hazel Round Metal Coffee Table In Powder Coating Finish", "₹19,763")
    ,@(37.0, "Peter Square Metal Coffee Table In Powder Coating Finish", "₹19,763")
    ,@(38.0, "Doug Square Metal Coffee Table In Powder Coating Finish", "₹19,763")
    ,@(39.0, "Angel Round Metal Coffee Table In Powder Coating Finish", "₹20,656")
    ,@(40.0, "Blane Square Solid Wood Coffee Table In Teak Finish", "₹19,999")
    ,@(41.0, "Blane Square Solid Wood Coffee Table In Antique Grey Finish", "₹19,499")
    ,@(42.0, "Blane Square Solid Wood Coffee Table In Antique Grey Finish", "₹19,499")
    ,@(43.0, "Blane Square Solid Wood Coffee Table In Antique Grey Finish", "₹19,499")
    ,@(44.0, "Nashville Round Solid Wood Coffee Table In Antique Grey Finish", "₹18,998")
    ,@(45.0, "Nashville Round Solid Wood Coffee Table In Antique Grey Finish", "₹18,998")
    ,@(46.0, "Nashville Round Solid Wood Coffee Table In Antique Grey Finish", "₹18,998")
    ,@(47.0, "Nashville Round Solid Wood Coffee Table In Antique Grey Finish", "₹18,998")
    ,@(48.0, "Blane Square Solid Wood Coffee Table In Teak Finish", "₹19,499")
    ,@(49.0, "Blane Square Solid Wood Coffee Table In Walnut Finish", "₹19,499")
    ,@(50.0, "Nashville Round Solid Wood Coffee Table In Teak Finish", "₹18,998")
    ,@(51.0, "Nashville Round Solid Wood Coffee Table In Teak Finish", "₹18,998")
    ,@(52.0, "Nashville Round Solid Wood Coffee Table In Teak Finish", "₹18,998")
    ,@(53.0, "Nashville Round Solid Wood Coffee Table In Teak Finish", "₹18,998")
    ,@(54.0, "Irish Rectangular Solid Wood Coffee Table In Walnut Finish", "₹13,999")
    ,@(55.0, "Irish Rectangular Solid Wood Coffee Table In Teak Finish", "₹13,999")
    ,@(56.0, "Irish Rectangular Solid Wood Coffee Table In Teak Finish", "₹13,999")
    ,@(57.0, "Irish Rectangular Solid Wood Coffee Table In Teak Finish", "₹13,999")
    ,@(58.0, "Irish Rectangular Solid Wood Coffee Table In Teak Finish", "₹13,999")
    ,@(59.0, "Milan Square Solid Wood Coffee Table In Walnut Finish", "₹13,498")
    ,@(60.0, "Montreal Square Solid Wood Coffee Table In Walnut Finish", "₹11,498")
)

$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

